# Insert a new data row for "Femacal de La Calera - Zanahoria" before the
# existing row 297. This shifts the existing rows 297-358 down to 298-359
# (preserving all their original values/formatting) and fills the newly
# inserted row 297 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 297, pushing rows 297-358 down to 298-359.
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row 297 with the new record's values.
$ws.Cells.Item(297, 1).Value  = 3
$ws.Cells.Item(297, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(297, 3).Value  = "Coquimbo"
$ws.Cells.Item(297, 4).Value  = 44711
$ws.Cells.Item(297, 5).Value  = 5
$ws.Cells.Item(297, 6).Value  = 100114013
$ws.Cells.Item(297, 7).Value  = "Zanahoria"
$ws.Cells.Item(297, 8).Value  = "Sin especificar"
$ws.Cells.Item(297, 9).Value  = "Primera"
$ws.Cells.Item(297, 10).Value = 570
$ws.Cells.Item(297, 11).Value = 7000
$ws.Cells.Item(297, 12).Value = 7500
$ws.Cells.Item(297, 13).Value = 7254
$ws.Cells.Item(297, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(297, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(297, 16).Value = 363
$ws.Cells.Item(297, 17).Value = 20
$ws.Cells.Item(297, 18).Value = "Hortaliza"
